$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Objective" heading -> "Summary"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Objective", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Summary", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Replace the old rambling "Objective" body paragraph with the new
#    Summary body text (single run, no line breaks).
# ---------------------------------------------------------------------------
$pSummaryBody = $d.Paragraphs.Item(3)
$pSummaryBody.Range.Text = "A passionate and driven student with a strong desire to contribute to the vibrant world of game development and UI design. My academic background in computer science and a keen understanding of game mechanics and user experience principles makes me a confident and adaptable individual eager to learn and contribute to innovative projects. I am committed to developing engaging and user-friendly applications, and I am eager to leverage my skills to create impactful solutions."

# ---------------------------------------------------------------------------
# 3. Append the new sections: Experience, Education, Skills, Projects.
#    Each new paragraph is produced via InsertParagraphAfter() on the
#    range of the previously-inserted paragraph, then its style/text set.
# ---------------------------------------------------------------------------

# --- Experience heading ---
$pSummaryBody.Range.InsertParagraphAfter()
$pExpHeading = $d.Paragraphs.Item(4)
$pExpHeading.Style = "Heading 1"
$pExpHeading.Range.Text = "Experience"

# --- Experience body ---
$pExpHeading.Range.InsertParagraphAfter()
$pExpBody = $d.Paragraphs.Item(5)
$pExpBody.Style = "Normal"
$pExpBody.Range.Text = "Here's a polished and professional version of the work experience section, incorporating the requested requirements:`v`v`v`v**AI Intern, Blue Silicon Infotech**`v`v`v`vHighly motivated and results-oriented AI intern with a proven track record of developing and deploying cutting-edge AI solutions. Developed expertise in natural language processing (NLP), machine learning (ML), and deep learning techniques. Successfully collaborated with diverse teams to design, develop, and implement AI models for various applications, including customer support chatbots and personalized recommendations. Demonstrated proficiency in data preprocessing, model training, and evaluation. Adept at collaborating effectively with stakeholders to achieve impactful outcomes. Opportunity to contribute to innovative AI projects and contribute to the growth of Blue Silicon Infotech's AI capabilities."

# --- Education heading ---
$pExpBody.Range.InsertParagraphAfter()
$pEduHeading = $d.Paragraphs.Item(6)
$pEduHeading.Style = "Heading 1"
$pEduHeading.Range.Text = "Education"

# --- Education body ---
$pEduHeading.Range.InsertParagraphAfter()
$pEduBody = $d.Paragraphs.Item(7)
$pEduBody.Style = "Normal"
$pEduBody.Range.Text = "B.E. CSE with 7.1 CGPA"

# --- Skills heading ---
$pEduBody.Range.InsertParagraphAfter()
$pSkillsHeading = $d.Paragraphs.Item(8)
$pSkillsHeading.Style = "Heading 1"
$pSkillsHeading.Range.Text = "Skills"

# --- Skills body ---
$pSkillsHeading.Range.InsertParagraphAfter()
$pSkillsBody = $d.Paragraphs.Item(9)
$pSkillsBody.Style = "Normal"
$pSkillsBody.Range.Text = "java, react, C#, node,java and python"

# --- Projects heading ---
$pSkillsBody.Range.InsertParagraphAfter()
$pProjHeading = $d.Paragraphs.Item(10)
$pProjHeading.Style = "Heading 1"
$pProjHeading.Range.Text = "Projects"

# --- Projects body ---
$pProjHeading.Range.InsertParagraphAfter()
$pProjBody = $d.Paragraphs.Item(11)
$pProjBody.Style = "Normal"
$pProjBody.Range.Text = "Here's the improved content for the project:`v`v**QR Scanner`v`v**`v`v**Prediction Pro**`v`v**Simple Purchase Order Manager**`v`v**Simple Purchase Order Manager**"

Write-Host ("Final paragraph count=" + $d.Paragraphs.Count)
